# Updates the cryptos price list: refreshed Price / Volume(1h) figures for
# row 2 to row 51, plus a NEARProtocol/Dai row swap (rows 26-27).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.577.13"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.041.96"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'202.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'629.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.45%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").Value = "'0.211"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").Value = "3.040.19"
$ws.Range("E10").Value = "  +4.31%  "
$ws.Range("D11").Value = "'0.439"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "'5.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").Value = "3.599.94"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("D15").Value = "'29.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.93%  "
$ws.Range("D16").Value = "76.427.12"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "3.027.38"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").Value = "'13.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").Value = "'9.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.67%  "
$ws.Range("D21").Value = "'376.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "'2.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").Value = "'73.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("D25").Value = "3.184.09"
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'4.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.35%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "'9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("E29").Value = "  +4.48%  "
$ws.Range("D30").Value = "'0.996"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'8.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.68%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").Value = "'514.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("E34").Value = "  +8.20%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'20.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.71%  "
$ws.Range("D37").Value = "'163.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  +11.05%  "
$ws.Range("D39").Value = "'20.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'0.108"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").Value = "'188.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "'5.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.82%  "
$ws.Range("E45").Value = "  +6.69%  "
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("E48").Value = "  +11.22%  "
$ws.Range("D49").Value = "'2.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.71%  "
$ws.Range("E50").Value = "  +6.48%  "
$ws.Range("D51").Value = "'3.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.06%  "
